$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the old summary rows (16:20) so nothing stale is left behind
# once the new layout (rows 18-23) is written.
$ws.Range("A16:D20").Clear()

# New column E: sequential index 1..14 next to the raw data rows.
for ($i = 0; $i -lt 14; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 5).Value = $i + 1
}

# Column A labels first (matches author's original shared-string order)
$ws.Range("A18").Value = "AVG"
$ws.Range("A19").Value = "STD"
$ws.Range("A20").Value = "Err AVG"
$ws.Range("A22").Value = "Disper"
$ws.Range("A23").Value = "%"

# Row 18: AVG
$ws.Range("B18").Formula = "=AVERAGE(B2:B15)"
$ws.Range("C18").Formula = "=AVERAGE(C2:C15)"

# Row 19: STD
$ws.Range("B19").Formula = "=STDEV.P(B2:B15)"
$ws.Range("C19").Formula = "=STDEV.P(C2:C15)"

# Row 20: Err AVG
$ws.Range("B20").Formula = "=B19/SQRT(14)"
$ws.Range("C20").Formula = "=C19/SQRT(14)"

# Row 21 intentionally left blank (only carries the shaded style below)

# Row 22: Disper
$ws.Range("B22").Formula = "=B19/B18"
$ws.Range("C22").Formula = "=C19/C18"

# Row 23: %
$ws.Range("B23").Formula = "=B22*100"
$ws.Range("C23").Formula = "=C22*100"

# Column D labels (reuses shared strings where the text already exists)
$ws.Range("D18").Value = "Width"
$ws.Range("D19").Value = "STD"
$ws.Range("D20").Value = "ERR Width"
$ws.Range("D22").Value = "Disper"
$ws.Range("D23").Value = "%"

# Shading: column A + D of the summary block get the blue / green fills.
$blue = 15773696   # RGB(0, 176, 240) -> FF00B0F0
$green = 5296274   # RGB(146, 208, 80) -> FF92D050

$ws.Range("A18:A23").Interior.Color = $blue
$ws.Range("D18:D23").Interior.Color = $green

# Selection as left by the author after the edit.
$ws.Range("B6:C8").Select()

Write-Host "done"
